# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, [string]$value) {
    # Force the cell to remain a plain text value (matches the original
    # inlineStr cells) instead of letting Excel auto-convert numeric-
    # looking strings (e.g. "536.92") into real numbers, then restore
    # the default "Normal" style/number format so no stray styles are left.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-CellText $ws.Range("D2") "59.214.57"
Set-CellText $ws.Range("E2") "  -0.28%  "
Set-CellText $ws.Range("D3") "2.526.13"
Set-CellText $ws.Range("E3") "  +0.29%  "
Set-CellText $ws.Range("E4") "  -0.10%  "
Set-CellText $ws.Range("D5") "536.92"
Set-CellText $ws.Range("E5") "  -0.96%  "
Set-CellText $ws.Range("D6") "136.90"
Set-CellText $ws.Range("E6") "  -2.25%  "
Set-CellText $ws.Range("E7") "  +0.06%  "
Set-CellText $ws.Range("D8") "0.569"
Set-CellText $ws.Range("E8") "  +0.66%  "
Set-CellText $ws.Range("D9") "2.525.72"
Set-CellText $ws.Range("E9") "  +0.10%  "
Set-CellText $ws.Range("E10") "  -0.24%  "
Set-CellText $ws.Range("E11") "  -2.51%  "
Set-CellText $ws.Range("D12") "5.31"
Set-CellText $ws.Range("E12") "  -1.93%  "
Set-CellText $ws.Range("D13") "0.349"
Set-CellText $ws.Range("E13") "  -0.65%  "
Set-CellText $ws.Range("D14") "2.974.29"
Set-CellText $ws.Range("E14") "  -0.06%  "
Set-CellText $ws.Range("D15") "23.09"
Set-CellText $ws.Range("E15") "  -1.28%  "
Set-CellText $ws.Range("D16") "59.149.08"
Set-CellText $ws.Range("E16") "  -0.30%  "
Set-CellText $ws.Range("E17") "  -1.51%  "
Set-CellText $ws.Range("D18") "2.531.02"
Set-CellText $ws.Range("E18") "  +0.48%  "
Set-CellText $ws.Range("D19") "11.15"
Set-CellText $ws.Range("E19") "  +0.44%  "
Set-CellText $ws.Range("E20") "  -0.02%  "
Set-CellText $ws.Range("D21") "324.01"
Set-CellText $ws.Range("E21") "  -0.45%  "
Set-CellText $ws.Range("E22") "  +0.07%  "
Set-CellText $ws.Range("E23") "  +1.31%  "
Set-CellText $ws.Range("D24") "65.51"
Set-CellText $ws.Range("E24") "  +3.40%  "
Set-CellText $ws.Range("E25") "  -0.08%  "
Set-CellText $ws.Range("E26") "  -1.84%  "
Set-CellText $ws.Range("E27") "  -0.03%  "
Set-CellText $ws.Range("E28") "  -2.66%  "
Set-CellText $ws.Range("E29") "  -1.41%  "
Set-CellText $ws.Range("E30") "  -0.46%  "
Set-CellText $ws.Range("B31") "Monero"
Set-CellText $ws.Range("C31") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText $ws.Range("D31") "171.98"
Set-CellText $ws.Range("E31") "  +4.04%  "
Set-CellText $ws.Range("B32") "PancakeSwap"
Set-CellText $ws.Range("C32") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText $ws.Range("D32") "1.76"
Set-CellText $ws.Range("E32") "  -1.85%  "
Set-CellText $ws.Range("E33") "  +5.94%  "
Set-CellText $ws.Range("E34") "  +0.01%  "
Set-CellText $ws.Range("D35") "1.42"
Set-CellText $ws.Range("E35") "  +0.40%  "
Set-CellText $ws.Range("D36") "18.39"
Set-CellText $ws.Range("E36") "  -0.88%  "
Set-CellText $ws.Range("D37") "4.11"
Set-CellText $ws.Range("E37") "  -1.87%  "
Set-CellText $ws.Range("E38") "  -2.94%  "
Set-CellText $ws.Range("E39") "  -0.69%  "
Set-CellText $ws.Range("D40") "0.814"
Set-CellText $ws.Range("E40") "  +0.71%  "
Set-CellText $ws.Range("E41") "  -2.00%  "
Set-CellText $ws.Range("D42") "285.49"
Set-CellText $ws.Range("E42") "  +1.56%  "
Set-CellText $ws.Range("E43") "  -2.29%  "
Set-CellText $ws.Range("E44") "  +0.15%  "
Set-CellText $ws.Range("E45") "  +1.64%  "
Set-CellText $ws.Range("D46") "131.20"
Set-CellText $ws.Range("E46") "  +3.75%  "
Set-CellText $ws.Range("D47") "10.87"
Set-CellText $ws.Range("E47") "  -0.04%  "
Set-CellText $ws.Range("E48") "  -1.61%  "
Set-CellText $ws.Range("E49") "  -1.07%  "
Set-CellText $ws.Range("E50") "  -1.47%  "
Set-CellText $ws.Range("D51") "17.43"
Set-CellText $ws.Range("E51") "  -2.41%  "
